# SDS_learning_diary_template.docx - add the final "18.5.2025" entry,
# AI declaration and acknowledgement sections at the end of the
# learning diary, after the "16.5.2025" entry's last paragraph
# ("I coded only the back end as the react part was similar to react
# module.").

$d = $word.ActiveDocument

# Grab the very last paragraph in the document (the "16.5.2025" entry's
# closing sentence) and collapse to its end so we can append after it.
$tailRange = $d.Paragraphs.Last.Range
$tailRange.Collapse(0)

# Build the new block of paragraphs. Each "`r" starts a brand new
# paragraph (mirroring Word's own paragraph-mark character), so an empty
# string between two "`r" markers yields a blank paragraph, same as the
# blank lines already used throughout the diary for spacing.
$newContent = (
    "`r" + `
    "`r" + `
    "18.5.2025`r" + `
    "Recorded a video to showcase my project. Probably fix some things yet. `r" + `
    "The project is still work in progress.`r" + `
    "Changes will be made but I will not update this document as this is only for the course. `r" + `
    "`r" + `
    "`r" + `
    "AI declaration:`r" + `
    "`r" + `
    "Some css is made by ChatGPT 4.0 or Deepseek. It is marked as such in those files. Otherwise no other Ais were used as 18-05-2025. `r" + `
    "`r" + `
    "Acknowledgement:`r" + `
    "Similarities might be seen with this project and Advanced web programming course and the module videos as everything is mainly learned from Advanced web programming course / those modules in this course. "
)

$tailRange.InsertAfter($newContent)
